# Generate Report for Handback
# Reflects a.md having been handed back (target xliff received, in sync
# with en-US source) for both the zh-cn and de-de localization targets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/360cfff5ea70736d7978372665bcb73c21c7fab0/e2e/a.md"

$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$zhHandbackDate = "2016-09-01 18:43:43"
$deHandbackDate = "2016-09-01 18:43:50"

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-locale status cells and widen the two
# status columns (E = zh-cn, F = de-de) so the longer text is legible.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet: mark both rows (a.md, b.md) as handed back, record the
# returned target file + handback datetime, and link the "Latest Target
# File" cell back to a.md (the only file with a produced target so far).
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("J3").Value = $zhHandbackFile

$wsZh.Range("K2").Value = $zhHandbackDate
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet: same treatment, with its own handback file/date.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("J3").Value = $deHandbackFile

$wsDe.Range("K2").Value = $deHandbackDate
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15
